$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "82÷6=13, 4" "65÷6=10, 5"
Replace-Text "53÷4=13, 1" "60÷6=10, 0"
Replace-Text "60÷7=8, 4" "37÷3=12, 1"
Replace-Text "76÷2=38, 0" "58÷4=14, 2"
Replace-Text "61÷2=30, 1" "94÷4=23, 2"
Replace-Text "97÷3=32, 1" "36÷6=6, 0"
Replace-Text "17÷7=2, 3" "33÷5=6, 3"
Replace-Text "89÷4=22, 1" "85÷6=14, 1"
Replace-Text "96÷6=16, 0" "74÷3=24, 2"
Replace-Text "51÷2=25, 1" "15÷3=5, 0"
Replace-Text "40÷4=10, 0" "52÷6=8, 4"
Replace-Text "90÷7=12, 6" "30÷7=4, 2"
Replace-Text "21÷2=10, 1" "36÷6=6, 0"
Replace-Text "67÷7=9, 4" "10÷5=2, 0"
Replace-Text "32÷2=16, 0" "62÷4=15, 2"
Replace-Text "54÷2=27, 0" "68÷6=11, 2"
Replace-Text "47÷6=7, 5" "60÷9=6, 6"
Replace-Text "70÷2=35, 0" "43÷7=6, 1"
Replace-Text "86÷7=12, 2" "83÷3=27, 2"
Replace-Text "25÷4=6, 1" "19÷2=9, 1"
Replace-Text "62÷9=6, 8" "66÷4=16, 2"
Replace-Text "30÷9=3, 3" "20÷3=6, 2"
Replace-Text "18÷2=9, 0" "41÷7=5, 6"
Replace-Text "85÷4=21, 1" "92÷3=30, 2"
Replace-Text "67÷3=22, 1" "22÷2=11, 0"
